$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.470.36"
$ws.Range("E2").Value = "  -1.23%  "
$ws.Range("D3").Value = "1.833.24"
$ws.Range("E3").Value = "  -1.39%  "
$ws.Range("E4").Value = "  -3.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.15%  "
$ws.Range("E6").Value = "  -2.66%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4309"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.38%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3711"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.89%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07268"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8697"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.23"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.92%  "
$ws.Range("D12").Value = "1.847.05"
$ws.Range("E12").Value = "  -0.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.700"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.377"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07110"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.007"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008938"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.91%  "
$ws.Range("E19").Value = "  -2.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.78%  "
$ws.Range("D21").Value = "27.474.69"
$ws.Range("E21").Value = "  -1.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.181"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.42%  "
$ws.Range("D24").Value = "2.059.43"
$ws.Range("E24").Value = "  -1.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.011"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.70"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.35%  "
$ws.Range("E27").Value = "  -1.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.159"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.306"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.60"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08877"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.212"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7710"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.514"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.915"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.005"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.126"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.67%  "
$ws.Range("E38").Value = "  -0.92%  "
$ws.Range("E39").Value = "  -0.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.194"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.880"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1684"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5096"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.721"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.68"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "106.65"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.91%  "
$ws.Range("E47").Value = "  +0.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06428"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.46%  "
$ws.Range("E49").Value = "  -2.94%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.678"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.838"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.89%  "
